# Commit: "Fruta / hortaliza, semanal"
# Insert two new weekly price records for Papa (Asterix / Rosara, cosecha,
# fecha 44931) above the existing row 628, pushing the rest of the table
# down by two rows (old 628..666 -> new 630..668).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 628.
$ws.Range("A628:A629").EntireRow.Insert()

# New row 628: Asterix, 1a (cosecha)
$ws.Range("A628").Value = 5
$ws.Range("B628").Value = "Macroferia Regional de Talca"
$ws.Range("C628").Value = "Maule"
$ws.Range("D628").Value = 44931
$ws.Range("E628").Value = 7
$ws.Range("F628").Value = 100114001
$ws.Range("G628").Value = "Papa"
$ws.Range("H628").Value = "Asterix"
$ws.Range("I628").Value = "1a (cosecha)"
$ws.Range("J628").Value = 1600
$ws.Range("K628").Value = 11000
$ws.Range("L628").Value = 11000
$ws.Range("M628").Value = 11000
$ws.Range("N628").Value = "`$/saco 25 kilos"
$ws.Range("O628").Value = "Región del Maule"
$ws.Range("P628").Value = 440
$ws.Range("Q628").Value = 25
$ws.Range("R628").Value = "Hortaliza"

# New row 629: Rosara, 1a (cosecha)
$ws.Range("A629").Value = 5
$ws.Range("B629").Value = "Macroferia Regional de Talca"
$ws.Range("C629").Value = "Maule"
$ws.Range("D629").Value = 44931
$ws.Range("E629").Value = 7
$ws.Range("F629").Value = 100114001
$ws.Range("G629").Value = "Papa"
$ws.Range("H629").Value = "Rosara"
$ws.Range("I629").Value = "1a (cosecha)"
$ws.Range("J629").Value = 1600
$ws.Range("K629").Value = 10000
$ws.Range("L629").Value = 10000
$ws.Range("M629").Value = 10000
$ws.Range("N629").Value = "`$/saco 25 kilos"
$ws.Range("O629").Value = "Región del Maule"
$ws.Range("P629").Value = 400
$ws.Range("Q629").Value = 25
$ws.Range("R629").Value = "Hortaliza"
